# Auto-generated edit script: updates cryptos list cell values
# per commit "Updated cryptos list on Thu Jul 20 23:56:31 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.820.47"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").Value = "1.892.59"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'0.7935"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.41%  "

# Row 6
$ws.Range("D6").Value = "'242.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.62%  "

# Row 8
$ws.Range("D8").Value = "'0.3161"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.57%  "

# Row 9
$ws.Range("D9").Value = "'25.36"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.92%  "

# Row 10
$ws.Range("E10").Value = "  +0.54%  "

# Row 11
$ws.Range("E11").Value = "  +0.51%  "

# Row 12
$ws.Range("D12").Value = "'0.7671"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.86%  "

# Row 13
$ws.Range("D13").Value = "1.891.42"
$ws.Range("E13").Value = "  +0.01%  "

# Row 14
$ws.Range("D14").Value = "'5.352"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.86%  "

# Row 15
$ws.Range("D15").Value = "'92.38"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.16%  "

# Row 16
$ws.Range("D16").Value = "29.842.08"
$ws.Range("E16").Value = "  -0.34%  "

# Row 17
$ws.Range("D17").Value = "'6.009"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.11%  "

# Row 18
$ws.Range("D18").Value = "'13.86"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.19%  "

# Row 19
$ws.Range("D19").Value = "'244.32"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.13%  "

# Row 20
$ws.Range("D20").Value = "'0.000007698"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.80%  "

# Row 21
$ws.Range("D21").Value = "'8.345"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +20.51%  "

# Row 22
$ws.Range("E22").Value = "  -0.07%  "

# Row 23
$ws.Range("D23").Value = "2.147.30"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("D25").Value = "'0.1633"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.88%  "

# Row 26
$ws.Range("D26").Value = "'9.362"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.80%  "

# Row 27
$ws.Range("D27").Value = "'166.06"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.00%  "

# Row 28
$ws.Range("D28").Value = "'18.69"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.65%  "

# Row 29
$ws.Range("D29").Value = "'2.053"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.52%  "

# Row 30
$ws.Range("E30").Value = "  +2.43%  "

# Row 31
$ws.Range("D31").Value = "'1.541"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.66%  "

# Row 32
$ws.Range("D32").Value = "'4.433"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.93%  "

# Row 33
$ws.Range("E33").Value = "  +1.27%  "

# Row 34
$ws.Range("D34").Value = "'4.045"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.53%  "

# Row 35
$ws.Range("D35").Value = "'1.262"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.59%  "

# Row 36
$ws.Range("D36").Value = "'0.7380"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.47%  "

# Row 37
$ws.Range("D37").Value = "'0.9992"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("D38").Value = "'2.632"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.08%  "

# Row 39
$ws.Range("D39").Value = "'0.01911"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("E40").Value = "  +0.12%  "

# Row 41
$ws.Range("D41").Value = "'0.4405"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.43%  "

# Row 42
$ws.Range("D42").Value = "'72.48"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.91%  "

# Row 43
$ws.Range("D43").Value = "'5.815"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.34%  "

# Row 44
$ws.Range("D44").Value = "'0.8414"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.27%  "

# Row 45
$ws.Range("E45").Value = "  -0.03%  "

# Row 46
$ws.Range("D46").Value = "1.034.06"
$ws.Range("E46").Value = "  +4.37%  "

# Row 47
$ws.Range("D47").Value = "'103.15"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.56%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.977"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.72%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.869"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("D50").Value = "'7.423"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.87%  "

# Row 51
$ws.Range("D51").Value = "2.047.76"
$ws.Range("E51").Value = "  +0.41%  "

